# Applies the 江西-漫展信息 update (commit "456a3b4") to both the
# "展览" and "全部类型" worksheets, which carry identical data tables.
#
# Summary of the change:
#  - A brand-new event ("赣州·第三届半夏动漫展") is inserted as row 15,
#    pushing the former rows 15-17 down to rows 16-18 (dimension grows
#    from A1:I17 to A1:I18).
#  - A handful of "想去人数" (F column) attendance counters drift by a
#    few units on rows 2,4,5,7,9,11,12 (and again on the shifted rows).

function Set-TextValue {
    # Writes a plain-text value into a cell while stopping Excel's
    # automatic "looks like a date" (or similar) reinterpretation, then
    # restores the cell to the unstyled "Normal" style so no stray
    # number-format survives the round trip (matches source cells, which
    # carry no explicit style).
    param($ws, $addr, [string]$val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

function Set-RowA {
    # Column A holds the running index and uses the bold/bordered/centered
    # "s=1" style. Copy that formatting from a cell that already has it
    # before writing the new index so no duplicate style gets synthesized.
    param($ws, $row, $index)
    $srcAddr = "A" + ($row - 1)
    $dstAddr = "A" + $row
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
    $ws.Range($dstAddr).Value = $index
}

function Apply-SheetEdits {
    param($ws)

    # --- Minor attendance-count ("想去人数", column F) drifts ---
    $ws.Range("F2").Value = 1569
    $ws.Range("F4").Value = 1026
    $ws.Range("F5").Value = 24
    $ws.Range("F7").Value = 2643
    $ws.Range("F9").Value = 1667
    $ws.Range("F11").Value = 66
    $ws.Range("F12").Value = 551

    # --- Row 18 (new): carries what used to be row 17's content,
    #     unchanged. Build its column-A formatting/value first, from the
    #     row above, then fill in the rest of the row. ---
    Set-RowA $ws 18 17
    Set-TextValue $ws "B18" "2024-04-13"
    $ws.Range("C18").Value = "南昌·第二届漫拥动漫嘉年华mini"
    $ws.Range("D18").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Range("E18").Value = "2024.04.13 10:00-04.14 18:00"
    $ws.Range("F18").Value = 8
    $ws.Range("G18").Value = 39.9
    $ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=82210"
    $ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png"

    # --- Row 17: now holds what used to be row 16's content (the
    #     "南昌·原X穹X崩only" event), with its attendance count nudged
    #     from 78 to 79. A17/B17 (index 16, date 2024-04-13) stay as-is. ---
    $ws.Range("C17").Value = "南昌·原X穹X崩only"
    $ws.Range("D17").Value = "丰和北大道299号 新吉花园酒店"
    $ws.Range("E17").Value = "2024.04.13 10:00-04.13 17:00"
    $ws.Range("F17").Value = 79
    $ws.Range("G17").Value = 65
    $ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=80807"
    $ws.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg"

    # --- Row 16: now holds what used to be row 15's content (the
    #     "赣州·赣次元·归来国风动漫节" event), with its date moved to
    #     2024-04-04 and attendance nudged from 64 to 67. A16 (index 15)
    #     stays as-is. ---
    Set-TextValue $ws "B16" "2024-04-04"
    $ws.Range("C16").Value = "赣州·赣次元·归来国风动漫节"
    $ws.Range("D16").Value = "客家大道568号文清外国语学校旁 赣州市文清外国语学校国际交流中心"
    $ws.Range("E16").Value = "2024.04.04 10:00-04.04 17:00"
    $ws.Range("F16").Value = 67
    $ws.Range("G16").Value = 40
    $ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=82125"
    $ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202402/8RNepTak1709022774421.jpeg"

    # --- Row 15: replaced with the brand-new event. A15/B15 (index 14,
    #     date 2024-04-04) stay as-is. ---
    $ws.Range("C15").Value = "赣州·第三届半夏动漫展"
    $ws.Range("D15").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
    $ws.Range("E15").Value = "2024.04.04 10:00-04.06 17:00"
    $ws.Range("F15").Value = 23
    $ws.Range("G15").Value = 36.6
    $ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=82235"
    $ws.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202403/4DWZWYGm1709278879159.jpeg"
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Apply-SheetEdits $ws
}
